$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 0.6142369984511111
$ws.Range("R2").Value = 5.52813298606
$ws.Range("S2").Value = 0.005922909240051363
$ws.Range("T2").Value = 0.005922909240051365

$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 0.2121874288572222
$ws.Range("R3").Value = 1.909686859715
$ws.Range("S3").Value = 0.00204606183959987
$ws.Range("T3").Value = 0.002046061839599871

$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 0.06953667366111112
$ws.Range("R4").Value = 0.6258300629500001
$ws.Range("S4").Value = 0.0006705219776542203
$ws.Range("T4").Value = 0.0006705219776542205

$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 57.50560574757067
$ws.Range("R5").Value = 517.550451728136
$ws.Range("S5").Value = 0.5545098789162999
$ws.Range("T5").Value = 0.5545098789163

$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("Q6").Value = 19.86524201443934
$ws.Range("S6").Value = 0.1915547675895075
$ws.Range("T6").Value = 0.1915547675895076

$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 6.510106930446668
$ws.Range("R7").Value = 58.59096237402001
$ws.Range("S7").Value = 0.06277507312209554
$ws.Range("T7").Value = 0.06277507312209556

$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.752937333333333
$ws.Range("N8").Value = 11.258812
$ws.Range("O8").Value = 0.6855621274031838
$ws.Range("P8").Value = 0.6855621274031838
$ws.Range("Q8").Value = 12.97657395075644
$ws.Range("R8").Value = 116.789165556808
$ws.Range("S8").Value = 0.1251293392468324
$ws.Range("T8").Value = 0.1251293392468324

$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.2368266084628361
$ws.Range("P9").Value = 0.2368266084628362
$ws.Range("Q9").Value = 4.482741790106889
$ws.Range("R9").Value = 40.34467611096201
$ws.Range("S9").Value = 0.0432257790337287
$ws.Range("T9").Value = 0.04322577903372871

$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 0.4248633333333334
$ws.Range("N10").Value = 1.27459
$ws.Range("O10").Value = 0.07761126413398003
$ws.Range("P10").Value = 0.07761126413398005
$ws.Range("Q10").Value = 1.469054762784445
$ws.Range("R10").Value = 13.22149286506
$ws.Range("S10").Value = 0.01416566903423027
$ws.Range("T10").Value = 0.01416566903423027
